$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mae")

$data = @(
    @{Row=5; J=84.211330083565429; K=236.75963888888879},
    @{Row=6; J=84.641448467966583; K=218.08497222222221},
    @{Row=7; J=84.239004178272964; K=215.69088888888879},
    @{Row=8; J=83.298426183844001; K=231.21124999999989},
    @{Row=9; J=81.161559888579404; K=243.80224999999999},
    @{Row=10; J=82.456427576601669; K=208.3794722222222},
    @{Row=11; J=85.929136490250684; K=222.09616666666659},
    @{Row=12; J=80.713426183844035; K=254.3535277777778},
    @{Row=13; J=79.992562674094714; K=276.0599444444444},
    @{Row=14; J=86.018962395543198; K=216.13900000000001},
    @{Row=15; J=77.790480501392764; K=289.13411111111111},
    @{Row=16; J=87.079860724233981; K=196.1835833333333},
    @{Row=17; J=85.728830083565484; K=211.1418333333333},
    @{Row=18; J=83.499707520891349; K=207.15077777777771},
    @{Row=19; J=81.553955431754858; K=232.8100833333333},
    @{Row=20; J=81.184394150417816; K=247.53949999999989},
    @{Row=21; J=87.361378830083567; K=214.29624999999999},
    @{Row=22; J=83.945431754874662; K=247.0813055555555},
    @{Row=23; J=85.912465181058494; K=209.3821388888889},
    @{Row=24; J=83.985076601671324; K=221.24941666666669},
    @{Row=25; J=85.065313370473532; K=203.49197222222219},
    @{Row=26; J=84.341267409470746; K=217.5622222222222},
    @{Row=27; J=80.737381615598878; K=260.94008333333329},
    @{Row=28; J=84.611608635097497; K=229.2571388888889},
    @{Row=29; J=84.272444289693581; K=234.03530555555551},
    @{Row=30; J=84.459993036211714; K=220.39105555555551},
    @{Row=31; J=79.952005571030654; K=264.38722222222219},
    @{Row=32; J=81.004213091922011; K=218.9157222222222},
    @{Row=33; J=84.973635097493045; K=199.01019444444441},
    @{Row=34; J=84.822367688022311; K=230.27483333333331},
    @{Row=35; J=85.490041782729804; K=220.35172222222221},
    @{Row=36; J=81.164679665738163; K=244.88677777777781},
    @{Row=37; J=82.983272980501368; K=213.05055555555549},
    @{Row=38; J=84.37545961002786; K=228.0026388888889},
    @{Row=39; J=81.887576601671313; K=251.6039999999999},
    @{Row=40; J=84.849784122562667; K=219.03522222222219},
    @{Row=41; J=83.412764623955439; K=221.97311111111111},
    @{Row=42; J=85.350090529247908; K=199.0651666666667},
    @{Row=43; J=82.092688022284122; K=246.18861111111119},
    @{Row=44; J=84.253140668523642; K=240.04116666666661},
    @{Row=45; J=85.119136490250725; K=228.43741666666659},
    @{Row=46; J=83.520222841225646; K=210.82222222222219},
    @{Row=47; J=79.662889972144853; K=255.65613888888879},
    @{Row=48; J=82.973809192200562; K=231.96905555555551},
    @{Row=49; J=84.377367688022289; K=231.58827777777779},
    @{Row=50; J=85.085731197771565; K=224.00155555555551},
    @{Row=51; J=82.008850974930368; K=270.76558333333332},
    @{Row=52; J=87.941357938718653; K=189.07472222222219},
    @{Row=53; J=84.578934540389966; K=236.17058333333341},
    @{Row=54; J=81.476497214484667; K=246.71838888888891}

)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 10).Value = $item.J
    $ws.Cells.Item($item.Row, 11).Value = $item.K
}

# Update the sheet view: scroll the window so row 28 is at the top, and
# move/select the active cell to F60 (matches the saved view state).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F60").Select()
